# Daily attendance processing - 2025-10-16 23:40:03
# Reorders the "Recorded By" (column G) contributor lists on the
# "Session Analysis Results" sheet so that "System"/"system" no longer
# sorts first, and "dnasr281@gmail.com" is listed ahead of
# "admin@admin.com".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Session Analysis Results")

$colG = 7

# Literal old -> new text replacements observed for the "Recorded By" column.
$replacements = @{
    "System, system, backup@backdoor.com" = "backup@backdoor.com, System, system"
    "System, dnasr281@gmail.com"          = "dnasr281@gmail.com, System"
    "System, backup@backdoor.com"         = "backup@backdoor.com, System"
    "admin@admin.com, dnasr281@gmail.com" = "dnasr281@gmail.com, admin@admin.com"
}

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

for ($row = 2; $row -le $lastRow; $row++) {
    $cell = $ws.Cells.Item($row, $colG)
    $value = $cell.Value2
    if ($null -ne $value -and $replacements.ContainsKey($value)) {
        $cell.Value = $replacements[$value]
    }
}
